# Fix power of two error in length of SF12 symbol
#
# Column F converts the raw receive-window timing (column E, in ms) into a
# number of SF12 LoRa symbol periods. The symbol period constant used was
# 16.384 ms, but the correct SF12 symbol period is 32.768 ms, so every
# formula in F2:F46 needs its divisor corrected. F3 is the anchor of a
# shared formula spanning F3:F46, and F2 is its own standalone formula, so
# updating those two ranges updates every dependent cell (and the chart
# that plots column F) consistently.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Formula = "=E2/32.768"
$ws.Range("F3:F46").Formula = "=E3/32.768"

# Recalculate so every cached <v> reflects the corrected divisor.
$excel.CalculateFullRebuild()

# Restore the cursor/selection position left by the author after making
# the edit (scrolled down, clicked the first empty row past the table).
$ws.Range("A49").Select()
